# Update "想去人数" (interest count) values in column F for the sheets
# "展览" and "全部类型", mirroring the regenerated site data dump.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value
$updates = @{
    2  = 1091
    5  = 4659
    8  = 1393
    9  = 921
    11 = 1138
    13 = 624
    15 = 34
    16 = 12
    18 = 24
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
